# Apply scheduled profit-recalculation updates across all Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 279.2857
$ws.Range("I11").Value = 279.2857
$ws.Range("K11").Value = 279.2857
$ws.Range("M11").Value = -139.2857
$ws.Range("H17").Value = 756.5769
$ws.Range("J17").Value = 756.5769
$ws.Range("L17").Value = 2269.7307
$ws.Range("N17").Value = -2605.7307
$ws.Range("H38").Value = 2243.7334
$ws.Range("I38").Value = 150.54546
$ws.Range("K38").Value = 451.63638
$ws.Range("M38").Value = -79.63637999999997
$ws.Range("H40").Value = 4177.8
$ws.Range("I40").Value = 2281.2856
$ws.Range("K40").Value = 2281.2856
$ws.Range("M40").Value = -2106.2856
$ws.Range("H92").Value = 365.9091
$ws.Range("I92").Value = 365.9091
$ws.Range("K92").Value = 365.9091
$ws.Range("M92").Value = 882.0908999999999
$ws.Range("H107").Value = 639.1905
$ws.Range("I107").Value = 662.15
$ws.Range("J107").Value = 180
$ws.Range("K107").Value = 662.15
$ws.Range("L107").Value = 180
$ws.Range("M107").Value = 1257.85
$ws.Range("N107").Value = -4020
$ws.Range("H138").Value = 24392100
$ws.Range("I138").Value = 1699.7894
$ws.Range("K138").Value = 5099.3682
$ws.Range("M138").Value = 40.63180000000011

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5251.914
$ws.Range("I32").Value = 5111.136
$ws.Range("K32").Value = 5111.136
$ws.Range("M32").Value = -4824.136
$ws.Range("H37").Value = 18000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null
$ws.Range("H45").Value = 5807.75
$ws.Range("I45").Value = 2448.8333
$ws.Range("K45").Value = 2448.8333
$ws.Range("M45").Value = -2071.8333
$ws.Range("H61").Value = 13895767
$ws.Range("I61").Value = 15630915
$ws.Range("K61").Value = 15630915
$ws.Range("M61").Value = -15630703
$ws.Range("H74").Value = 2812.0967
$ws.Range("I74").Value = 2292
$ws.Range("K74").Value = 2292
$ws.Range("M74").Value = -1418
$ws.Range("H77").Value = 2812.0967
$ws.Range("I77").Value = 2292
$ws.Range("K77").Value = 11460
$ws.Range("M77").Value = -7092
$ws.Range("H80").Value = 43000
$ws.Range("J80").Value = 43000
$ws.Range("L80").Value = 43000
$ws.Range("N80").Value = -44996
$ws.Range("H83").Value = 43000
$ws.Range("J83").Value = 43000
$ws.Range("L83").Value = 129000
$ws.Range("N83").Value = -138984
$ws.Range("H97").Value = 2580.6
$ws.Range("I97").Value = 3230
$ws.Range("K97").Value = 3230
$ws.Range("M97").Value = -2734
$ws.Range("H102").Value = 3740.4827
$ws.Range("I102").Value = 3263.6956
$ws.Range("J102").Value = 5568.1665
$ws.Range("K102").Value = 3263.6956
$ws.Range("L102").Value = 5568.1665
$ws.Range("M102").Value = -1641.6956
$ws.Range("N102").Value = -8812.166499999999
$ws.Range("H122").Value = 1511.4546
$ws.Range("I122").Value = 1141.4706
$ws.Range("K122").Value = 3424.4118
$ws.Range("M122").Value = -974.4118000000003
$ws.Range("H132").Value = 4311.561
$ws.Range("I132").Value = 3789.1
$ws.Range("J132").Value = 5736.4546
$ws.Range("K132").Value = 11367.3
$ws.Range("L132").Value = 17209.3638
$ws.Range("M132").Value = -8837.299999999999
$ws.Range("N132").Value = -22269.3638
$ws.Range("H136").Value = 13895767
$ws.Range("I136").Value = 15630915
$ws.Range("K136").Value = 46892745
$ws.Range("M136").Value = -46890195

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 8749.75
$ws.Range("I22").Value = 8499.5
$ws.Range("K22").Value = 8499.5
$ws.Range("M22").Value = -8326.5
$ws.Range("H94").Value = 4538
$ws.Range("I94").Value = 4538
$ws.Range("K94").Value = 4538
$ws.Range("M94").Value = -4087
$ws.Range("H95").Value = 52374.75
$ws.Range("J95").Value = 52374.75
$ws.Range("L95").Value = 52374.75
$ws.Range("N95").Value = -57866.75
$ws.Range("H134").Value = 5294.3335
$ws.Range("I134").Value = 5294.3335
$ws.Range("K134").Value = 15883.0005
$ws.Range("M134").Value = -13348.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4996
$ws.Range("I6").Value = 4994.6665
$ws.Range("K6").Value = 4994.6665
$ws.Range("M6").Value = -4881.6665
$ws.Range("H7").Value = 42.5
$ws.Range("I7").Value = 46.11111
$ws.Range("K7").Value = 46.11111
$ws.Range("M7").Value = 66.88889
$ws.Range("H107").Value = 554.08
$ws.Range("I107").Value = 496.33334
$ws.Range("K107").Value = 496.33334
$ws.Range("M107").Value = 1423.66666
$ws.Range("H134").Value = 6036.9614
$ws.Range("I134").Value = 4741.467
$ws.Range("K134").Value = 14224.401
$ws.Range("M134").Value = -11689.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 137815.86
$ws.Range("J37").Value = 137815.86
$ws.Range("L37").Value = 413447.58
$ws.Range("N37").Value = -413671.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15033.111
$ws.Range("I70").Value = 11449.5
$ws.Range("K70").Value = 11449.5
$ws.Range("M70").Value = -11179.5
$ws.Range("H73").Value = 15033.111
$ws.Range("I73").Value = 11449.5
$ws.Range("K73").Value = 11449.5
$ws.Range("M73").Value = -10513.5
$ws.Range("H102").Value = 4517.75
$ws.Range("I102").Value = 4028.4
$ws.Range("K102").Value = 4028.4
$ws.Range("M102").Value = -2406.4
$ws.Range("H132").Value = 5231.1
$ws.Range("I132").Value = 4039
$ws.Range("K132").Value = 12117
$ws.Range("M132").Value = -9587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1632.5555
$ws.Range("I22").Value = 3197
$ws.Range("J22").Value = 1185.5714
$ws.Range("K22").Value = 3197
$ws.Range("L22").Value = 1185.5714
$ws.Range("M22").Value = -2902
$ws.Range("N22").Value = -1775.5714
$ws.Range("H27").Value = 1632.5555
$ws.Range("I27").Value = 3197
$ws.Range("J27").Value = 1185.5714
$ws.Range("K27").Value = 3197
$ws.Range("L27").Value = 1185.5714
$ws.Range("M27").Value = -3090
$ws.Range("N27").Value = -1399.5714
$ws.Range("H40").Value = 2850.5
$ws.Range("I40").Value = 2459.4119
$ws.Range("K40").Value = 2459.4119
$ws.Range("M40").Value = -2323.4119
$ws.Range("H122").Value = 3474
$ws.Range("I122").Value = 3474
$ws.Range("K122").Value = 10422
$ws.Range("M122").Value = -7972
$ws.Range("H132").Value = 8954.583000000001
$ws.Range("I132").Value = 8949.093999999999
$ws.Range("K132").Value = 26847.282
$ws.Range("M132").Value = -24317.282
$ws.Range("H136").Value = 4416.522
$ws.Range("I136").Value = 3740.5
$ws.Range("K136").Value = 11221.5
$ws.Range("M136").Value = -8671.5
$ws.Range("H137").Value = 85100
$ws.Range("J137").Value = 85100
$ws.Range("L137").Value = 85100
$ws.Range("N137").Value = -95300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10923.571
$ws.Range("I5").Value = 1500
$ws.Range("J5").Value = 12494.167
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 12494.167
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = -12718.167
$ws.Range("H122").Value = 4899.4287
$ws.Range("I122").Value = 2649
$ws.Range("K122").Value = 7947
$ws.Range("M122").Value = -5497
$ws.Range("H126").Value = 1801.1765
$ws.Range("I126").Value = 1801.1765
$ws.Range("K126").Value = 5403.529500000001
$ws.Range("M126").Value = -2933.529500000001
$ws.Range("H132").Value = 4516.5
$ws.Range("I132").Value = 4018.15
$ws.Range("K132").Value = 12054.45
$ws.Range("M132").Value = -9524.450000000001
$ws.Range("H136").Value = 4202.6665
$ws.Range("I136").Value = 3696.138
$ws.Range("K136").Value = 11088.414
$ws.Range("M136").Value = -8538.414000000001
